# Update BOC USD rates (auto)
# A new rate was published for 2026-01-02 (17:53:08) on top of the one
# already captured at 17:43:08, so:
#   - "All Published Values" gains a new row (row 5) with the fresh reading
#   - "Daily Summary" day-averages table now counts 2 publishes for that day

$wb = $excel.ActiveWorkbook

# --- Sheet "All Published Values": append the new published-rate row ---
$wsAll = $wb.Worksheets.Item("All Published Values")

# Every column on this sheet stores plain text (even numeric/date-looking
# values), so prefix with an apostrophe to stop Excel from re-interpreting
# them as numbers/dates, then strip the resulting "Text" number format so
# the cells end up unstyled just like the rest of the sheet.
$newRow = 5
$wsAll.Cells.Item($newRow, 1).Value  = "'2026-01-02"
$wsAll.Cells.Item($newRow, 2).Value  = "'2026-01-02 17:53:08"
$wsAll.Cells.Item($newRow, 3).Value  = "'697.85"
$wsAll.Cells.Item($newRow, 4).Value  = "'697.85"
$wsAll.Cells.Item($newRow, 5).Value  = "'700.79"
$wsAll.Cells.Item($newRow, 6).Value  = "'700.79"
$wsAll.Cells.Item($newRow, 7).Value  = "'702.88"
$wsAll.Cells.Item($newRow, 8).Value  = "'2026/01/02 17:53:08"
$wsAll.Cells.Item($newRow, 9).Value  = "'2026-01-02 09:56:30"
$wsAll.Cells.Item($newRow, 10).Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
$wsAll.Range("A5:J5").ClearFormats()

# --- Sheet "Daily Summary": 2026-01-02 now has 2 publishes instead of 1 ---
# (this column is stored as text too, so apply the same apostrophe trick)
$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Cells.Item(4, 2).Value = "'2"
$wsSummary.Cells.Item(4, 2).ClearFormats()
